# Actualización automática del inventario: se agrega un nuevo producto
# (fila 78) a la hoja de inventario, igual que las filas existentes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 78

$ws.Range("A$row").Value = "8RLBN2"
$ws.Range("B$row").Value = "Led con base 3 volt"
# Columna C (Descripción) se deja vacía, igual que en la fila anterior.
$ws.Range("D$row").Value = 2500
$ws.Range("E$row").Value = 7000
$ws.Range("F$row").Value = 197
$ws.Range("G$row").Value = 53
$ws.Range("H$row").Formula = "=(E$row-D$row)*G$row"
$ws.Range("I$row").Formula = "=D$row*F$row"
$ws.Range("J$row").Value = 492500
